$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values.
# Leading apostrophe forces these numeric-looking strings to be stored as text,
# matching the source data which stores all values as plain text (inline strings).
$ws.Range("D2").Value = "'289.86"
$ws.Range("E2").Value = "'-3.81%"
$ws.Range("D3").Value = "'30.51"
$ws.Range("E3").Value = "'-5.69%"
$ws.Range("D4").Value = "'4.945"
$ws.Range("E4").Value = "'-1.08%"
$ws.Range("E5").Value = "'-5.83%"
$ws.Range("D6").Value = "'1.821"
$ws.Range("E6").Value = "'-6.20%"
$ws.Range("D7").Value = "'7.619"
$ws.Range("E7").Value = "'-2.65%"
$ws.Range("D8").Value = "'3.718"
$ws.Range("E8").Value = "'-1.69%"
$ws.Range("D9").Value = "'0.8962"
$ws.Range("E9").Value = "'-2.22%"
$ws.Range("D10").Value = "'0.1652"
$ws.Range("E10").Value = "'-5.69%"
$ws.Range("D11").Value = "'0.07728"
$ws.Range("E11").Value = "'-1.03%"
$ws.Range("D12").Value = "'0.07946"
$ws.Range("E12").Value = "'-7.07%"
$ws.Range("D13").Value = "'0.03040"
$ws.Range("E13").Value = "'-3.28%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("E15").Value = "'-0.70%"
$ws.Range("D16").Value = "'0.005735"
$ws.Range("E16").Value = "'0.16%"
$ws.Range("D18").Value = "'3.464"
$ws.Range("E18").Value = "'0.04%"
$ws.Range("E20").Value = "'-0.67%"
$ws.Range("D21").Value = "'0.1278"
$ws.Range("E21").Value = "'-3.55%"
$ws.Range("D22").Value = "'4.014"
$ws.Range("E22").Value = "'-6.11%"
$ws.Range("D23").Value = "'0.2389"
$ws.Range("E23").Value = "'20.01%"
$ws.Range("D24").Value = "'0.04506"
$ws.Range("E24").Value = "'-0.17%"
$ws.Range("E25").Value = "'-0.48%"
$ws.Range("D26").Value = "'0.004618"
$ws.Range("E26").Value = "'5.16%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'4.04%"
$ws.Range("D39").Value = "'0.01563"
$ws.Range("E39").Value = "'-8.32%"
$ws.Range("D40").Value = "'0.04349"
$ws.Range("E40").Value = "'-7.05%"
$ws.Range("D41").Value = "'0.007321"
$ws.Range("E41").Value = "'-2.15%"
$ws.Range("D43").Value = "'0.1301"
$ws.Range("E43").Value = "'-3.53%"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("E44").Value = "'-11.54%"
$ws.Range("D45").Value = "'0.009505"
$ws.Range("E45").Value = "'-9.92%"
$ws.Range("D46").Value = "'0.00005985"
$ws.Range("E46").Value = "'-4.35%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("D48").Value = "'2.255"
$ws.Range("E48").Value = "'174.73%"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E51").Value = "'0.05%"

Write-Output "Updated crypto price/volume cells."
